# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (positioned between the existing
# "2021-Q4" and "总计" sheets) populated with per-fund holding data, and
# updates the "总计" (totals) sheet with a new summary row for 2022-Q1
# (the old 2021-Q4 summary row shifts down).

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Add($null, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

# NOTE: worksheet references are positional, so "总计" must be re-fetched
# *after* the insert above shifts its index (it was index 2, now index 3).
$totalSheet = $wb.Worksheets.Item("总计")

# Row data for the new sheet: index, fund code, fund name, fund scale,
# total stock position, position ratio, holding value (亿元), position rank
$rows = @(
    @(0, "014269", "嘉实北交所精选两年定期混合A",           "5.00", "87.43", "6.23", "0.3115", 3),
    @(1, "014283", "华夏北交所创新中小企业精选两年定开混合", "3.96", "52.55", "5.34", "0.2115", 4),
    @(2, "008962", "建信科技创新混合A",                     "4.34", "84.82", "4.73", "0.2053", 3),
    @(3, "014273", "广发北交所精选两年定开混合A",           "4.55", "52.69", "4.32", "0.1966", 6),
    @(4, "014294", "南方北交所精选两年定开混合",             "4.63", "33.00", "1.80", "0.0833", 4),
    @(5, "014663", "富国创新发展两年定期开放混合A",         "2.62", "37.11", "2.34", "0.0613", 2),
    @(6, "014270", "嘉实北交所精选两年定期混合C",           "0.64", "87.43", "6.23", "0.0399", 3),
    @(7, "014274", "广发北交所精选两年定开混合C",           "0.92", "52.69", "4.32", "0.0397", 6),
    @(8, "008963", "建信科技创新混合C",                     "0.26", "84.82", "4.73", "0.0123", 3),
    @(9, "014664", "富国创新发展两年定期开放混合C",         "0.32", "37.11", "2.34", "0.0075", 2)
)

# Copy the header row (B1:H1) formatting + text straight from "2021-Q4" --
# both sheets share an identical header ("基金代码", "基金名称", ...).
$q4Sheet.Range("B1:H1").Copy($q1Sheet.Range("B1"))

# Copy the data-row template (A2:H2) down for every one of the 10 rows so
# each row starts with the right styles (bordered/bold index column A,
# plain text columns B:G, plain numeric column H).
for ($i = 2; $i -le 11; $i++) {
    $q4Sheet.Range("A2:H2").Copy($q1Sheet.Range("A" + $i + ":H" + $i))
}

# Force columns B:G to Text so numeric-looking strings (fund codes with
# leading zeros, "5.00", "0.3115", ...) are stored as text, not numbers.
$q1Sheet.Range("B2:G11").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q1Sheet.Cells.Item($r, 1).Value = $row[0]
    $q1Sheet.Cells.Item($r, 2).Value = $row[1]
    $q1Sheet.Cells.Item($r, 3).Value = $row[2]
    $q1Sheet.Cells.Item($r, 4).Value = $row[3]
    $q1Sheet.Cells.Item($r, 5).Value = $row[4]
    $q1Sheet.Cells.Item($r, 6).Value = $row[5]
    $q1Sheet.Cells.Item($r, 7).Value = $row[6]
    $q1Sheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing 2021-Q4 summary row down
#    to row 3 (and bump its index 0 -> 1), then write the new 2022-Q1
#    summary row into row 2.
# ---------------------------------------------------------------------
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))
$totalSheet.Range("A3").Value = 1

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 10
$totalSheet.Range("D2").Value = 1.17
